$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(47, "Create Country", "PASSED", "chrome", "13.12.22")
  ,@(48, "Add Document Types", "PASSED", "chrome", "16.12.22")
  ,@(49, "Add Document Types", "PASSED", "chrome", "17.12.22")
  ,@(50, "Edit Document Types", "FAILED", "chrome", "17.12.22")
  ,@(51, "Edit Document Types", "FAILED", "chrome", "17.12.22")
  ,@(52, "Edit Document Types", "FAILED", "chrome", "17.12.22")
  ,@(53, "Edit Document Types", "FAILED", "chrome", "17.12.22")
  ,@(54, "Edit Document Types", "FAILED", "chrome", "18.12.22")
  ,@(55, "Edit Document Types", "FAILED", "chrome", "18.12.22")
  ,@(56, "Add Document Types", "PASSED", "chrome", "18.12.22")
  ,@(57, "Edit Document Types", "FAILED", "chrome", "18.12.22")
  ,@(58, "Add Document Types", "PASSED", "chrome", "18.12.22")
  ,@(59, "Edit Document Types", "FAILED", "chrome", "18.12.22")
  ,@(60, "Edit Document Types", "PASSED", "chrome", "18.12.22")
  ,@(61, "Add Document Types", "PASSED", "chrome", "18.12.22")
  ,@(62, "Edit Document Types", "FAILED", "chrome", "18.12.22")
  ,@(63, "Add Document Types", "PASSED", "chrome", "19.12.22")
  ,@(64, "Edit Document Types", "FAILED", "chrome", "19.12.22")
  ,@(65, "Delete Document Types", "FAILED", "chrome", "19.12.22")
  ,@(66, "Edit Document Types", "FAILED", "chrome", "19.12.22")
  ,@(67, "Edit Document Types", "FAILED", "chrome", "19.12.22")
  ,@(68, "Edit Document Types", "PASSED", "chrome", "19.12.22")
  ,@(69, "Delete Document Types", "PASSED", "chrome", "19.12.22")
  ,@(70, "Add Document Types", "PASSED", "chrome", "19.12.22")
  ,@(71, "Edit Document Types", "FAILED", "chrome", "19.12.22")
  ,@(72, "Add Document Types", "FAILED", "chrome", "19.12.22")
  ,@(73, "Edit Document Types", "FAILED", "chrome", "19.12.22")
  ,@(74, "Add Document Types", "PASSED", "chrome", "20.12.22")
  ,@(75, "Add Document Types", "PASSED", "chrome", "20.12.22")
  ,@(76, "Edit Document Types", "FAILED", "chrome", "20.12.22")
  ,@(77, "Delete Document Types", "PASSED", "chrome", "20.12.22")
  ,@(78, "Negative Delete Document Types", "PASSED", "chrome", "20.12.22")
  ,@(79, "Login with valid username and password", "PASSED", "chrome", "22.12.22")
  ,@(80, "Login with valid username and password", "PASSED", "chrome", "22.12.22")
  ,@(81, "Add Document Types", "PASSED", "chrome", "17.01.23")
  ,@(82, "Edit Document Types", "PASSED", "chrome", "17.01.23")
  ,@(83, "Delete Document Types", "PASSED", "chrome", "17.01.23")
  ,@(84, "Create A Citizenship", "PASSED", "chrome", "17.01.23")
  ,@(85, "Create An Existant Citizenship", "PASSED", "chrome", "17.01.23")
  ,@(86, "Update the Citizenship", "PASSED", "chrome", "17.01.23")
  ,@(87, "Delete the Citizenship", "PASSED", "chrome", "17.01.23")
  ,@(88, "Search and delete an unavailable Citizenship", "PASSED", "chrome", "17.01.23")
  ,@(89, "Create Country", "PASSED", "chrome", "17.01.23")
  ,@(90, "Create a Attestations", "FAILED", "chrome", "17.01.23")
  ,@(91, "Edit a Attestations", "FAILED", "chrome", "17.01.23")
  ,@(92, "Delete a Attestations", "FAILED", "chrome", "17.01.23")
  ,@(93, "Add Document Types", "PASSED", "chrome", "17.01.23")
  ,@(94, "Add Document Types", "FAILED", "chrome", "17.01.23")
  ,@(95, "Add Document Types", "FAILED", "chrome", "17.01.23")
  ,@(96, "Add Document Types", "PASSED", "chrome", "17.01.23")
  ,@(97, "Add Document Types", "PASSED", "chrome", "24.01.23")
  ,@(98, "Add School Locations", "PASSED", "chrome", "24.01.23")
  ,@(99, "Add School Locations", "PASSED", "chrome", "24.01.23")
  ,@(100, "Add School Locations", "PASSED", "chrome", "24.01.23")
  ,@(101, "Add School Locations", "PASSED", "chrome", "24.01.23")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}
